$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6123.75
$ws.Range("I62").Value = 3165
$ws.Range("K62").Value = 3165
$ws.Range("M62").Value = -2541

$ws.Range("H65").Value = 6123.75
$ws.Range("I65").Value = 3165
$ws.Range("K65").Value = 15825
$ws.Range("M65").Value = -12705

$ws.Range("H121").Value = 1765.0303
$ws.Range("J121").Value = 1802.6875
$ws.Range("L121").Value = 5408.0625
$ws.Range("N121").Value = -8902.0625

$ws.Range("H129").Value = 897.44116
$ws.Range("J129").Value = 945.6774
$ws.Range("L129").Value = 2837.0322
$ws.Range("N129").Value = -12837.0322

$ws.Range("H132").Value = 125186.695
$ws.Range("I132").Value = 192741.05
$ws.Range("J132").Value = 6966.5835
$ws.Range("K132").Value = 578223.1499999999
$ws.Range("L132").Value = 20899.7505
$ws.Range("M132").Value = -575693.1499999999
$ws.Range("N132").Value = -25959.7505

$ws.Range("H135").Value = 1180.2307
$ws.Range("I135").Value = 1040.4
$ws.Range("J135").Value = 1646.3334
$ws.Range("K135").Value = 9363.6
$ws.Range("L135").Value = 14817.0006
$ws.Range("M135").Value = -6828.6
$ws.Range("N135").Value = -19887.0006

$ws.Range("H137").Value = 3938.8718
$ws.Range("I137").Value = 3232.303
$ws.Range("J137").Value = 7825
$ws.Range("K137").Value = 9696.909
$ws.Range("L137").Value = 23475
$ws.Range("M137").Value = -7146.909
$ws.Range("N137").Value = -28575

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1948.1578
$ws.Range("I61").Value = 1308.3572
$ws.Range("J61").Value = 3739.6
$ws.Range("K61").Value = 1308.3572
$ws.Range("L61").Value = 3739.6
$ws.Range("M61").Value = -1096.3572
$ws.Range("N61").Value = -4163.6

$ws.Range("H74").Value = 3902.7878
$ws.Range("I74").Value = 4056.72
$ws.Range("K74").Value = 4056.72
$ws.Range("M74").Value = -3182.72

$ws.Range("H77").Value = 3902.7878
$ws.Range("I77").Value = 4056.72
$ws.Range("K77").Value = 20283.6
$ws.Range("M77").Value = -15915.6

$ws.Range("H82").Value = 48300
$ws.Range("J82").Value = 48300
$ws.Range("L82").Value = 48300
$ws.Range("N82").Value = -49022

$ws.Range("H85").Value = 48300
$ws.Range("J85").Value = 48300
$ws.Range("L85").Value = 48300
$ws.Range("N85").Value = -50796

$ws.Range("H136").Value = 1948.1578
$ws.Range("I136").Value = 1308.3572
$ws.Range("J136").Value = 3739.6
$ws.Range("K136").Value = 3925.0716
$ws.Range("L136").Value = 11218.8
$ws.Range("M136").Value = -1375.0716
$ws.Range("N136").Value = -16318.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3141.0889
$ws.Range("I134").Value = 1673
$ws.Range("J134").Value = 8279.4
$ws.Range("K134").Value = 5019
$ws.Range("L134").Value = 24838.2
$ws.Range("M134").Value = -2484
$ws.Range("N134").Value = -29908.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 6851.409
$ws.Range("I134").Value = 7934.8667
$ws.Range("J134").Value = 4529.7144
$ws.Range("K134").Value = 23804.6001
$ws.Range("L134").Value = 13589.1432
$ws.Range("M134").Value = -21269.6001
$ws.Range("N134").Value = -18659.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2494.1667
$ws.Range("I122").Value = 449.46155
$ws.Range("J122").Value = 3516.5193
$ws.Range("K122").Value = 4045.15395
$ws.Range("L122").Value = 31648.6737
$ws.Range("M122").Value = -1595.15395
$ws.Range("N122").Value = -36548.6737

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 20173
$ws.Range("J34").Value = 20173
$ws.Range("L34").Value = 20173
$ws.Range("N34").Value = -20709

$ws.Range("H76").Value = 20173
$ws.Range("J76").Value = 20173
$ws.Range("L76").Value = 20173
$ws.Range("N76").Value = -20803

$ws.Range("H79").Value = 20173
$ws.Range("J79").Value = 20173
$ws.Range("L79").Value = 20173
$ws.Range("N79").Value = -22357

$ws.Range("H132").Value = 3801.75
$ws.Range("I132").Value = 2381
$ws.Range("J132").Value = 5222.5
$ws.Range("K132").Value = 7143
$ws.Range("L132").Value = 15667.5
$ws.Range("M132").Value = -4613
$ws.Range("N132").Value = -20727.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4221.0215
$ws.Range("I40").Value = 3935.2163
$ws.Range("J40").Value = 5278.5
$ws.Range("K40").Value = 3935.2163
$ws.Range("L40").Value = 5278.5
$ws.Range("M40").Value = -3799.2163
$ws.Range("N40").Value = -5550.5

$ws.Range("H81").Value = 80999.336
$ws.Range("J81").Value = 80999.336
$ws.Range("L81").Value = 80999.336
$ws.Range("N81").Value = -82995.336

$ws.Range("H84").Value = 80999.336
$ws.Range("J84").Value = 80999.336
$ws.Range("L84").Value = 242998.008
$ws.Range("N84").Value = -252982.008

$ws.Range("H132").Value = 6847.88
$ws.Range("I132").Value = 2401.2727
$ws.Range("J132").Value = 10341.643
$ws.Range("K132").Value = 7203.8181
$ws.Range("L132").Value = 31024.929
$ws.Range("M132").Value = -4673.8181
$ws.Range("N132").Value = -36084.929

$ws.Range("H136").Value = 3627.1892
$ws.Range("I136").Value = 1371.7142
$ws.Range("J136").Value = 6587.5
$ws.Range("K136").Value = 4115.142599999999
$ws.Range("L136").Value = 19762.5
$ws.Range("M136").Value = -1565.142599999999
$ws.Range("N136").Value = -24862.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8929840
$ws.Range("I81").Value = 11906036
$ws.Range("J81").Value = 1250
$ws.Range("K81").Value = 23812072
$ws.Range("L81").Value = 2500
$ws.Range("M81").Value = -23811011
$ws.Range("N81").Value = -4622

$ws.Range("H84").Value = 8929840
$ws.Range("I84").Value = 11906036
$ws.Range("J84").Value = 1250
$ws.Range("K84").Value = 119060360
$ws.Range("L84").Value = 12500
$ws.Range("M84").Value = -119055056
$ws.Range("N84").Value = -23108

$ws.Range("H132").Value = 11496439
$ws.Range("I132").Value = 976.8182
$ws.Range("J132").Value = 18521444
$ws.Range("K132").Value = 2930.4546
$ws.Range("L132").Value = 55564332
$ws.Range("M132").Value = -400.4546
$ws.Range("N132").Value = -55569392

$ws.Range("H136").Value = 3550.3076
$ws.Range("I136").Value = 813.625
$ws.Range("J136").Value = 7929
$ws.Range("K136").Value = 2440.875
$ws.Range("L136").Value = 23787
$ws.Range("M136").Value = 109.125
$ws.Range("N136").Value = -28887
